# "fixed export and fixing maps"
#
# The sheet previously carried a census-comparison table (1989 / 2002 / 2014
# columns plus a "(census results)" subtitle). The fix trims it down to just
# the most recent (2014) figure, renames the tab to the municipality name,
# and leaves the selection on A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused 1989 / 2002 columns (column D, "2014", slides left into B).
$ws.Range("B:C").EntireColumn.Delete()

# Drop the "(მოსახლეობის აღწერის შედეგებით)" / census-results subtitle row.
$ws.Range("A2").EntireRow.Delete()

# Give the sheet its real (municipality) name instead of the generic "1".
$ws.Name = "ბოლნისი"

# Match the saved cursor position from the source file.
$ws.Range("A2").Select()
